# Realestate Update resale numbers 2024-01-12 21:36
# Appends a new data row (row 52) to the CityResaleNum sheet, matching the
# pattern of the existing rows: text Date/Time/Weekday/Week columns (A-D)
# and numeric city-resale-number columns (E-T).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 52

# Text columns: Date, Time, Weekday, Week.
# "2024-01-12" and "01" look like a date / number to Excel's auto-detection,
# so a leading apostrophe is used to force them to be stored as plain text,
# matching the other rows in the sheet (t="inlineStr"/t="s" string cells).
$ws.Cells.Item($row, 1).Value = "'2024-01-12"
$ws.Cells.Item($row, 2).Value = "21:36:11"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "'01"

# Numeric columns: Beijing, Guangzhou, Suzhou, Hangzhou, Nanjing, Xi_an,
# Chengdu, Chongqing, Tianjin, Hefei, Fuzhou, Xiamen, Changsha, Shanghai,
# Shenzhen, Wuhan.
$ws.Cells.Item($row, 5).Value = 136658
$ws.Cells.Item($row, 6).Value = 142740
$ws.Cells.Item($row, 7).Value = 172107
$ws.Cells.Item($row, 8).Value = 148314
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119774
$ws.Cells.Item($row, 11).Value = 225059
$ws.Cells.Item($row, 12).Value = 253053
$ws.Cells.Item($row, 13).Value = 184974
$ws.Cells.Item($row, 14).Value = 110437
$ws.Cells.Item($row, 15).Value = 40947
$ws.Cells.Item($row, 16).Value = 30922
$ws.Cells.Item($row, 17).Value = 73110
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42803
$ws.Cells.Item($row, 20).Value = -1

$wb.Save()
